# Adicionado todas funcoes basicas positiva e negativa
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Planilha1: update the "Usuario" test value in row 2 (B2) to a new user
$ws1.Range("B2").Value = "usertest99"

# Planilha2: add a new (currently empty) row with word-wrap formatting,
# matching the added row 4 / cell B4 in the sheet
$ws2.Range("B4").WrapText = $true
$ws2.Range("B4").RowHeight = 17.25

# Restore selections to match the saved view state
[void]$ws1.Range("D5").Select()
[void]$ws2.Range("C3").Select()
